$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 1483.4073
$ws.Range("I32").Value = 1391.4286
$ws.Range("J32").Value = 1515.6
$ws.Range("K32").Value = 1391.4286
$ws.Range("L32").Value = 1515.6
$ws.Range("M32").Value = -1065.4286
$ws.Range("N32").Value = -2167.6

# Row 33
$ws.Range("H33").Value = 149.2
$ws.Range("I33").Value = 119.23077
$ws.Range("J33").Value = 204.85715
$ws.Range("K33").Value = 119.23077
$ws.Range("L33").Value = 204.85715
$ws.Range("M33").Value = 109.76923
$ws.Range("N33").Value = -662.85715

# Row 40
$ws.Range("H40").Value = 22799.8
$ws.Range("I40").Value = 100001
$ws.Range("J40").Value = 3499.5
$ws.Range("K40").Value = 100001
$ws.Range("L40").Value = 3499.5
$ws.Range("M40").Value = -99826
$ws.Range("N40").Value = -3849.5

# Row 43
$ws.Range("H43").Value = 3031.25
$ws.Range("I43").Value = 3133.3333
$ws.Range("K43").Value = 3133.3333
$ws.Range("M43").Value = -3064.3333

# Row 51
$ws.Range("H51").Value = 4779
$ws.Range("I51").Value = 3395
$ws.Range("J51").Value = 5125
$ws.Range("K51").Value = 3395
$ws.Range("L51").Value = 5125
$ws.Range("M51").Value = -2911
$ws.Range("N51").Value = -6093

# Row 76
$ws.Range("H76").Value = 2913.5334
$ws.Range("I76").Value = 2662.875
$ws.Range("K76").Value = 2662.875
$ws.Range("M76").Value = -2347.875

# Row 79
$ws.Range("H79").Value = 2913.5334
$ws.Range("I79").Value = 2662.875
$ws.Range("K79").Value = 2662.875
$ws.Range("M79").Value = -1570.875

# Row 112
$ws.Range("H112").Value = 1152.5405
$ws.Range("J112").Value = 1152.5405
$ws.Range("L112").Value = 3457.6215
$ws.Range("N112").Value = -5673.6215

# Row 136
$ws.Range("H136").Value = 40310
$ws.Range("J136").Value = 40310
$ws.Range("L136").Value = 40310
$ws.Range("N136").Value = -50510

# Row 141
$ws.Range("H141").Value = 2809
$ws.Range("I141").Value = 2134.65
$ws.Range("J141").Value = 9552.5
$ws.Range("K141").Value = 6403.950000000001
$ws.Range("L141").Value = 28657.5
$ws.Range("M141").Value = -1223.950000000001
$ws.Range("N141").Value = -39017.5

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 1020.1
$ws.Range("I5").Value = 877.8889
$ws.Range("K5").Value = 877.8889
$ws.Range("M5").Value = -765.8889

# Row 45
$ws.Range("H45").Value = 41668468
$ws.Range("I45").Value = 47620676
$ws.Range("K45").Value = 47620676
$ws.Range("M45").Value = -47620299

# Row 96
$ws.Range("H96").Value = 32229
$ws.Range("J96").Value = 32229
$ws.Range("L96").Value = 32229
$ws.Range("N96").Value = -37721

# Row 104
$ws.Range("H104").Value = 30643.2
$ws.Range("J104").Value = 30643.2
$ws.Range("L104").Value = 30643.2
$ws.Range("N104").Value = -37631.2

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 1020.1
$ws.Range("I4").Value = 877.8889
$ws.Range("K4").Value = 877.8889
$ws.Range("M4").Value = -762.8889

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 430.16666
$ws.Range("I7").Value = 296.2
$ws.Range("J7").Value = 1100
$ws.Range("K7").Value = 296.2
$ws.Range("L7").Value = 1100
$ws.Range("M7").Value = -183.2
$ws.Range("N7").Value = -1326

# Row 31
$ws.Range("H31").Value = 6276.75
$ws.Range("I31").Value = 2750.077
$ws.Range("J31").Value = 8270.087
$ws.Range("K31").Value = 2750.077
$ws.Range("L31").Value = 8270.087
$ws.Range("M31").Value = -2455.077
$ws.Range("N31").Value = -8860.087

# Row 34
$ws.Range("H34").Value = 6276.75
$ws.Range("I34").Value = 2750.077
$ws.Range("J34").Value = 8270.087
$ws.Range("K34").Value = 2750.077
$ws.Range("L34").Value = 8270.087
$ws.Range("M34").Value = -2548.077
$ws.Range("N34").Value = -8674.087

# Row 62
$ws.Range("H62").Value = 2659.9583
$ws.Range("I62").Value = 2523.1052
$ws.Range("K62").Value = 2523.1052
$ws.Range("M62").Value = -1899.1052

# Row 65
$ws.Range("H65").Value = 2659.9583
$ws.Range("I65").Value = 2523.1052
$ws.Range("K65").Value = 12615.526
$ws.Range("M65").Value = -9495.526

# Row 68
$ws.Range("H68").Value = 167500
$ws.Range("J68").Value = 167500
$ws.Range("L68").Value = 167500
$ws.Range("N68").Value = -168998

# Row 70
$ws.Range("H70").Value = 30209.143
$ws.Range("J70").Value = 30209.143
$ws.Range("L70").Value = 30209.143
$ws.Range("N70").Value = -30839.143

# Row 71
$ws.Range("H71").Value = 167500
$ws.Range("J71").Value = 167500
$ws.Range("L71").Value = 502500
$ws.Range("N71").Value = -509988

# Row 73
$ws.Range("H73").Value = 30209.143
$ws.Range("J73").Value = 30209.143
$ws.Range("L73").Value = 30209.143
$ws.Range("N73").Value = -32393.143

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 622.7273

# Row 71
$ws.Range("H71").Value = 622.7273

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 1205589.6
$ws.Range("I2").Value = 1721942
$ws.Range("J2").Value = 767.3333
$ws.Range("K2").Value = 1721942
$ws.Range("L2").Value = 767.3333
$ws.Range("M2").Value = -1721829
$ws.Range("N2").Value = -993.3333

# Row 102
$ws.Range("H102").Value = 1377.5834
$ws.Range("I102").Value = 815.2857
$ws.Range("J102").Value = 2164.8
$ws.Range("K102").Value = 815.2857
$ws.Range("L102").Value = 2164.8
$ws.Range("M102").Value = 806.7143
$ws.Range("N102").Value = -5408.8

# Row 132
$ws.Range("H132").Value = 3662.3
$ws.Range("I132").Value = 2923.8333
$ws.Range("J132").Value = 4770
$ws.Range("K132").Value = 8771.499899999999
$ws.Range("L132").Value = 14310
$ws.Range("M132").Value = -6241.499899999999
$ws.Range("N132").Value = -19370

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1349
$ws.Range("I22").Value = 1299
$ws.Range("J22").Value = 1399
$ws.Range("K22").Value = 1299
$ws.Range("L22").Value = 1399
$ws.Range("M22").Value = -1004
$ws.Range("N22").Value = -1989

# Row 27
$ws.Range("H27").Value = 1349
$ws.Range("I27").Value = 1299
$ws.Range("J27").Value = 1399
$ws.Range("K27").Value = 1299
$ws.Range("L27").Value = 1399
$ws.Range("M27").Value = -1192

# Row 46
$ws.Range("H46").Value = 1912.4375
$ws.Range("J46").Value = 2537.375
$ws.Range("L46").Value = 2537.375
$ws.Range("N46").Value = -2913.375

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 26429454
$ws.Range("I122").Value = 44048456
$ws.Range("J122").Value = 947.375
$ws.Range("K122").Value = 132145368
$ws.Range("L122").Value = 2842.125
$ws.Range("M122").Value = -132142918
$ws.Range("N122").Value = -7742.125
